$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: plain values
$ws.Range("F12").Value = 1305498328.7699957
$ws.Range("G12").Value = 1240524717

# Row 13: plain values
$ws.Range("F13").Value = 325268233.58999997
$ws.Range("G13").Value = 319819483.19999999

# Row 14: plain value
$ws.Range("G14").Value = 3406311.8

# Row 16: plain values
$ws.Range("F16").Value = -53616441.74000001
$ws.Range("G16").Value = 60834434.380000003

# Row 18: F18 becomes a formula (was a plain value); G18 formula stays, will recalc
$ws.Range("F18").Formula = "=SUM(F12:F17)"

# Row 19: plain value
$ws.Range("G19").Value = 379300000

# Row 21: F21 becomes a formula (was a plain value); G21 formula stays, will recalc
$ws.Range("F21").Formula = "=SUM(F18:F20)"

# Row 22: plain value
$ws.Range("G22").Value = 20015625

# Row 26: plain values (inputs to F25/G25 formulas)
$ws.Range("F26").Value = 1026703455.3810816
$ws.Range("G26").Value = 1029174575

# Force full recalculation so dependent formula cells (F23, G23, F25, G25,
# F28, G28, C29, G18, G21) pick up the new cached values.
$excel.CalculateFullRebuild()
